# Update header: report volume number and week-covering dates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 31   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/30/2024  Through  10/6/2024"

# Update weekly crime-complaint statistics table (rows 14-33)
$ws.Range("D14").Value = '0'
$ws.Range("E14").Value = '***.*'
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("N14").Value = -60
$ws.Range("C15").Value = '0'
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 7
$ws.Range("H15").Value = -57.142857142857
$ws.Range("J15").Value = 32
$ws.Range("K15").Value = -12.5
$ws.Range("L15").Value = -9.677419354838
$ws.Range("M15").Value = 115.384615384615
$ws.Range("N15").Value = 16.666666666666
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -32
$ws.Range("I16").Value = 244
$ws.Range("J16").Value = 246
$ws.Range("K16").Value = -0.813008130081
$ws.Range("L16").Value = 9.909909909909
$ws.Range("M16").Value = 57.419354838709
$ws.Range("N16").Value = -73.274917853231
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 28.571428571428
$ws.Range("F17").Value = 47
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = 51.612903225806
$ws.Range("I17").Value = 408
$ws.Range("J17").Value = 345
$ws.Range("K17").Value = 18.260869565217
$ws.Range("L17").Value = 63.855421686747
$ws.Range("M17").Value = 144.311377245509
$ws.Range("N17").Value = 42.160278745644
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -22.222222222222
$ws.Range("F18").Value = 33
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = 10
$ws.Range("I18").Value = 314
$ws.Range("J18").Value = 424
$ws.Range("K18").Value = -25.943396226415
$ws.Range("L18").Value = -23.414634146341
$ws.Range("M18").Value = -16.710875331565
$ws.Range("N18").Value = -83.063646170442
$ws.Range("C19").Value = 26
$ws.Range("D19").Value = 25
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 95
$ws.Range("G19").Value = 85
$ws.Range("H19").Value = 11.764705882352
$ws.Range("I19").Value = 1016
$ws.Range("J19").Value = 1057
$ws.Range("K19").Value = -3.878902554399
$ws.Range("L19").Value = -12.714776632302
$ws.Range("M19").Value = 104.016064257028
$ws.Range("N19").Value = -0.974658869395
$ws.Range("C20").Value = 20
$ws.Range("D20").Value = 17
$ws.Range("E20").Value = 17.647058823529
$ws.Range("F20").Value = 57
$ws.Range("G20").Value = 31
$ws.Range("H20").Value = 83.870967741935
$ws.Range("I20").Value = 431
$ws.Range("J20").Value = 385
$ws.Range("K20").Value = 11.948051948051
$ws.Range("L20").Value = 130.48128342246
$ws.Range("M20").Value = 93.273542600896
$ws.Range("N20").Value = -86.114690721649
$ws.Range("C21").Value = 66
$ws.Range("D21").Value = 64
$ws.Range("E21").Value = 3.125
$ws.Range("G21").Value = 210
$ws.Range("H21").Value = 20.47619047619
$ws.Range("I21").Value = 2445
$ws.Range("J21").Value = 2492
$ws.Range("K21").Value = -1.886035313001
$ws.Range("L21").Value = 7.851786501985
$ws.Range("M21").Value = 70.027816411682
$ws.Range("N21").Value = -66.126350789692
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -66.666666666666
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 23
$ws.Range("J22").Value = 27
$ws.Range("K22").Value = -14.814814814814
$ws.Range("L22").Value = -23.333333333333
$ws.Range("M22").Value = 666.666666666667
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = '0'
$ws.Range("E23").Value = '***.*'
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 20
$ws.Range("K23").Value = 66.666666666666
$ws.Range("L23").Value = 66.666666666666
$ws.Range("M23").Value = 100
$ws.Range("C24").Value = 53
$ws.Range("D24").Value = 45
$ws.Range("E24").Value = 17.777777777777
$ws.Range("F24").Value = 199
$ws.Range("G24").Value = 187
$ws.Range("H24").Value = 6.417112299465
$ws.Range("I24").Value = 2099
$ws.Range("J24").Value = 2178
$ws.Range("K24").Value = -3.627180899908
$ws.Range("L24").Value = -2.144522144522
$ws.Range("M24").Value = 76.832350463353
$ws.Range("C25").Value = 34
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = 36
$ws.Range("F25").Value = 146
$ws.Range("G25").Value = 122
$ws.Range("H25").Value = 19.672131147541
$ws.Range("I25").Value = 1394
$ws.Range("J25").Value = 1195
$ws.Range("K25").Value = 16.652719665272
$ws.Range("L25").Value = 23.581560283687
$ws.Range("C26").Value = 26
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 62.5
$ws.Range("F26").Value = 82
$ws.Range("G26").Value = 70
$ws.Range("H26").Value = 17.142857142857
$ws.Range("I26").Value = 746
$ws.Range("J26").Value = 681
$ws.Range("K26").Value = 9.544787077826
$ws.Range("L26").Value = 39.439252336448
$ws.Range("M26").Value = 33.691756272401
$ws.Range("C27").Value = '0'
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -57.142857142857
$ws.Range("J27").Value = 48
$ws.Range("K27").Value = -25
$ws.Range("L27").Value = -14.285714285714
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = '0'
$ws.Range("E28").Value = '***.*'
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 60
$ws.Range("I28").Value = 73
$ws.Range("K28").Value = -9.876543209876
$ws.Range("L28").Value = 14.0625
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = -100
$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 7
$ws.Range("K31").Value = -14.285714285714
$ws.Range("L31").Value = -45.454545454545
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = -100
$ws.Range("J33").Value = 9
$ws.Range("K33").Value = 44.444444444444
